$d = $word.ActiveDocument

$replacements = @(
    @("45×77=", "58×40="),
    @("97×99=", "71×20="),
    @("60×38=", "78×42="),
    @("15×62=", "40×47="),
    @("25×85=", "54×18="),
    @("81×58=", "33×67="),
    @("37×18=", "70×12="),
    @("53×88=", "28×17="),
    @("44×31=", "79×66="),
    @("38×43=", "64×59="),
    @("91×38=", "97×56="),
    @("27×58=", "88×25="),
    @("29×73=", "30×15="),
    @("90×50=", "70×30="),
    @("79×85=", "20×30="),
    @("46×38=", "71×49="),
    @("42×47=", "16×81="),
    @("13×73=", "54×38="),
    @("76×90=", "84×82="),
    @("93×31=", "96×23="),
    @("70×29=", "16×50="),
    @("94×11=", "58×14="),
    @("45×89=", "62×47="),
    @("65×29=", "26×90="),
    @("40×56=", "34×26=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
